$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update TestCases value (B2): 40 -> 44
$ws.Range("B2").Value = "44"

# Update Instance value (D2): Automation1 -> Automation3
$ws.Range("D2").Value = "Automation3"

# Update the active selection to D2 (was E2)
$ws.Range("D2").Select()
